# Rewrite the acceptance-test read-me:
#   - add a bold "Android Acceptance Test Documentation" title
#   - reword the "logged out" bullet
#   - split the LoginActivity/ViewAdsActivity sentence and mark the
#     CamelCase identifiers as spell-check exceptions (w:proofErr)
#   - add a "Steps to Run Acceptance test" lead-in and a new
#     hanging-indent paragraph with the actual steps
#   - reword the storage-permission bullet and its two sub-bullets
#   - drop the stray _GoBack bookmark from the middle of the doc and
#     leave it on the final (trailing) empty paragraph instead
#
# The body is rebuilt in one shot via Range.InsertXML: that lets us set
# paragraph/run XML (pPr, rPr, proofErr, bookmarks) exactly, the same
# way Word's own paste/XML-import pipeline would, rather than fighting
# Find/Replace across run boundaries. $d.Content.InsertXML() replaces
# everything between the start and end of the document but leaves the
# final <w:sectPr> (section properties) alone.
#
# Note: w:val="a3" below is this document's existing internal style id
# for the "List Paragraph" style (see word/styles.xml) -- Word resolves
# pPr/w:pStyle by that id, not by the display name, so every paragraph
# that should carry the List Paragraph style references "a3".

$d = $word.ActiveDocument

$bodyXml = @'
<w:p>
  <w:pPr>
    <w:ind w:left="3600"/>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>Android Acceptance Test Documentation</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>Before running acceptance test, logout of the android app:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:ind w:left="1440"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Acceptance Test starts with </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>LoginActivity</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">. If tester is logged in and start the test, because of the saved session of user, the application will start with </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>ViewAdsActivity</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>, which will fail the Acceptance Test.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:ind w:left="1440"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:ind w:left="1440"/>
  </w:pPr>
  <w:r>
    <w:t>Steps</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> to Run Acceptance test</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">.  </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:ind w:left="1440" w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t>Log out from the application</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> in your emulator/or mobile test</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> -&gt; Run </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>AllAcceptanceTest</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:ind w:left="1440"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>Must give storage permission to run the acceptance. Without it test</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> will stall</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve">: </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t>Tester should go to the Settin</w:t>
  </w:r>
  <w:r>
    <w:t>gs -&gt; Apps &amp; Notifications -&gt; UMBUY</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> -&gt; Permissions -&gt; Give permission to Storage.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t>Or if click Allow when the applica</w:t>
  </w:r>
  <w:r>
    <w:t>tion asks for permission when creating an ad</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> first time.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:ind w:left="1440"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Steps. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">After giving permission for the photos/gallery doing one of above options -&gt; Run </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>AllAcceptanceTest</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:ind w:left="1440"/>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$packageXml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $bodyXml + '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

# $d.Content is the document's main-story range (everything except the
# final section mark). Replacing it with a <w:body> that has no
# <w:sectPr> leaves the existing section properties (page size/margins)
# untouched, and Word re-lays the new paragraphs out from the top.
$d.Content.InsertXML($packageXml)
